$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.2115286666666667
$ws.Cells.Item(2, 8).Value = 0.634586
$ws.Cells.Item(2, 9).Value = 0.08153347995807345
$ws.Cells.Item(2, 10).Value = 0.08153347995807345
$ws.Cells.Item(2, 13).Value = 1.970755333333333
$ws.Cells.Item(2, 14).Value = 5.912266
$ws.Cells.Item(2, 15).Value = 0.178292931636603
$ws.Cells.Item(2, 16).Value = 0.178292931636603
$ws.Cells.Item(2, 17).Value = 0.4168712479862222
$ws.Cells.Item(2, 18).Value = 3.751841231876
$ws.Cells.Item(2, 19).Value = 0.01453684316825913
$ws.Cells.Item(2, 20).Value = 0.01453684316825913

# Row 3
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.2115286666666667
$ws.Cells.Item(3, 8).Value = 0.634586
$ws.Cells.Item(3, 9).Value = 0.08153347995807345
$ws.Cells.Item(3, 10).Value = 0.08153347995807345
$ws.Cells.Item(3, 15).Value = 0.5919535183013138
$ws.Cells.Item(3, 16).Value = 0.5919535183013139
$ws.Cells.Item(3, 17).Value = 1.384061609503778
$ws.Cells.Item(3, 18).Value = 12.456554485534
$ws.Cells.Item(3, 19).Value = 0.04826403032053123
$ws.Cells.Item(3, 20).Value = 0.04826403032053124

# Row 4
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.2115286666666667
$ws.Cells.Item(4, 8).Value = 0.634586
$ws.Cells.Item(4, 9).Value = 0.08153347995807345
$ws.Cells.Item(4, 10).Value = 0.08153347995807345
$ws.Cells.Item(4, 15).Value = 0.2297535500620831
$ws.Cells.Item(4, 16).Value = 0.2297535500620831
$ws.Cells.Item(4, 17).Value = 0.5371926316117779
$ws.Cells.Item(4, 18).Value = 4.834733684506
$ws.Cells.Item(4, 19).Value = 0.01873260646928308
$ws.Cells.Item(4, 20).Value = 0.01873260646928308

# Row 5
$ws.Cells.Item(5, 9).Value = 0.4735790235655714
$ws.Cells.Item(5, 10).Value = 0.4735790235655714
$ws.Cells.Item(5, 13).Value = 1.970755333333333
$ws.Cells.Item(5, 14).Value = 5.912266
$ws.Cells.Item(5, 15).Value = 0.178292931636603
$ws.Cells.Item(5, 16).Value = 0.178292931636603
$ws.Cells.Item(5, 17).Value = 2.421354745012666
$ws.Cells.Item(5, 18).Value = 21.792192705114
$ws.Cells.Item(5, 19).Value = 0.08443579247310562
$ws.Cells.Item(5, 20).Value = 0.08443579247310565

# Row 6
$ws.Cells.Item(6, 9).Value = 0.4735790235655714
$ws.Cells.Item(6, 10).Value = 0.4735790235655714
$ws.Cells.Item(6, 15).Value = 0.5919535183013138
$ws.Cells.Item(6, 16).Value = 0.5919535183013139
$ws.Cells.Item(6, 19).Value = 0.2803367691933408
$ws.Cells.Item(6, 20).Value = 0.2803367691933409

# Row 7
$ws.Cells.Item(7, 9).Value = 0.4735790235655714
$ws.Cells.Item(7, 10).Value = 0.4735790235655714
$ws.Cells.Item(7, 15).Value = 0.2297535500620831
$ws.Cells.Item(7, 16).Value = 0.2297535500620831
$ws.Cells.Item(7, 19).Value = 0.1088064618991249
$ws.Cells.Item(7, 20).Value = 0.108806461899125

# Row 8
$ws.Cells.Item(8, 9).Value = 0.4448874964763552
$ws.Cells.Item(8, 10).Value = 0.4448874964763552
$ws.Cells.Item(8, 13).Value = 1.970755333333333
$ws.Cells.Item(8, 14).Value = 5.912266
$ws.Cells.Item(8, 15).Value = 0.178292931636603
$ws.Cells.Item(8, 16).Value = 0.178292931636603
$ws.Cells.Item(8, 17).Value = 2.274658287183778
$ws.Cells.Item(8, 18).Value = 20.471924584654
$ws.Cells.Item(8, 19).Value = 0.07932029599523825
$ws.Cells.Item(8, 20).Value = 0.07932029599523827

# Row 9
$ws.Cells.Item(9, 9).Value = 0.4448874964763552
$ws.Cells.Item(9, 10).Value = 0.4448874964763552
$ws.Cells.Item(9, 15).Value = 0.5919535183013138
$ws.Cells.Item(9, 16).Value = 0.5919535183013139
$ws.Cells.Item(9, 19).Value = 0.2633527187874418
$ws.Cells.Item(9, 20).Value = 0.2633527187874418

# Row 10
$ws.Cells.Item(10, 9).Value = 0.4448874964763552
$ws.Cells.Item(10, 10).Value = 0.4448874964763552
$ws.Cells.Item(10, 15).Value = 0.2297535500620831
$ws.Cells.Item(10, 16).Value = 0.2297535500620831
$ws.Cells.Item(10, 19).Value = 0.1022144816936751
$ws.Cells.Item(10, 20).Value = 0.1022144816936751
